$wb = $excel.ActiveWorkbook

# --- Overview sheet: the c04d7826 handoff, which previously showed
# "Ready for handoff", now failed the handback transform. ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B7").Value = "Handback transform failed"
$overview.Range("C7").Value = "Handback transform failed"

# --- zh-cn sheet: record the Error Detail explaining the handback
# file name mismatch for the c04d7826 row (row 7). ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("L7").Value = "Handback file name: oeibrr3g.qrj is different with handoff file name: c04d7826-47bc-46f5-a399-ec30212a6fe9.c4590bbb13d2319e1020ceb20c358e10b6d926eb.zh-cn."

# --- de-de sheet: same Error Detail, for the de-de locale. ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("L7").Value = "Handback file name: oeibrr3g.qrj is different with handoff file name: c04d7826-47bc-46f5-a399-ec30212a6fe9.c4590bbb13d2319e1020ceb20c358e10b6d926eb.de-de."
